# The diff shows the survey records that were in rows 23-25 being
# cyclically rotated: the "Garnlav" record (originally row 23) moves to
# row 25, the "Spindelblomster (25 st)" record (originally row 24) moves
# to row 23, and the "Spindelblomster (3 st)" record (originally row 25)
# moves to row 24. Apply that by writing each record's field values into
# its new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 23 ----
$ws.Range('A23').Value = 111664006
$ws.Range('B23').Value = 96368
$ws.Range('D23').Value = 'LC'
$ws.Range('E23').Value = 221952
$ws.Range('F23').Value = 'Spindelblomster'
$ws.Range('G23').Value = 'Neottia cordata'
$ws.Range('H23').Value = '(L.) Rich.'
$ws.Range('I23').NumberFormat = "@"
$ws.Range('I23').Value = '25'
$ws.Range('J23').Value = 'stjälkar/strån/skott'
$ws.Range('K23').Value = 'överblommad'
$ws.Range('L23').Value = ''
$ws.Range('N23').Value = 'observerad'
$ws.Range('Q23').Value = 485633.0508789647
$ws.Range('R23').Value = 7005896.397059018
$ws.Range('S23').Value = 10
$ws.Range('AC23').Value = 'Minst 25 st. skott/stjälkar av spindelblomster varav minst ett skott som bär en överblommad blomstängel. Obs! Fyndplatsen ligger inom en avverkningsanmäld yta med beteckning A 32699-2023.'
$ws.Range('AJ23').Value = ''
$ws.Range('AK23').Value = ''
$ws.Range('AO23').Value = ''

# ---- Row 24 ----
$ws.Range('A24').Value = 111663810
$ws.Range('I24').NumberFormat = "@"
$ws.Range('I24').Value = '3'
$ws.Range('Q24').Value = 485613.4985160928
$ws.Range('R24').Value = 7005872.099109154
$ws.Range('S24').Value = 5
$ws.Range('AC24').Value = 'Minst 3 skott/stjälkar av spindelblomster. Obs! Fyndplatsen ligger inom en avverkningsanmäld yta med beteckning A 32699-2023.'

# ---- Row 25 ----
$ws.Range('A25').Value = 111663296
$ws.Range('B25').Value = 77515
$ws.Range('D25').Value = 'NT'
$ws.Range('E25').Value = 6425
$ws.Range('F25').Value = 'Garnlav'
$ws.Range('G25').Value = 'Alectoria sarmentosa'
$ws.Range('H25').Value = '(Ach.) Ach.'
$ws.Range('I25').Value = ''
$ws.Range('J25').Value = ''
$ws.Range('K25').Value = ''
$ws.Range('L25').Value = ''
$ws.Range('N25').Value = ''
$ws.Range('Q25').Value = 485609.6900141542
$ws.Range('R25').Value = 7005829.216201009
$ws.Range('AC25').Value = 'Enstaka bålar av garnlav på gran. Obs! Fyndplatsen ligger inom en avverkningsanmäld yta med beteckning A 32699-2023.'
$ws.Range('AJ25').Value = 'gran'
$ws.Range('AK25').Value = 'Picea abies'
$ws.Range('AO25').Value = 'Picea abies'

Write-Output "Rotated rows 23-25 per the source diff."
